$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-19 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-20 Saturday", 2) | Out-Null
$d.Content.Find.Execute("41×91=", $true, $false, $false, $false, $false, $true, 1, $false, "42×60=", 2) | Out-Null
$d.Content.Find.Execute("58×50=", $true, $false, $false, $false, $false, $true, 1, $false, "89×82=", 2) | Out-Null
$d.Content.Find.Execute("28×95=", $true, $false, $false, $false, $false, $true, 1, $false, "93×78=", 2) | Out-Null
$d.Content.Find.Execute("12×12=", $true, $false, $false, $false, $false, $true, 1, $false, "19×89=", 2) | Out-Null
$d.Content.Find.Execute("33×21=", $true, $false, $false, $false, $false, $true, 1, $false, "48×20=", 2) | Out-Null
$d.Content.Find.Execute("31×14=", $true, $false, $false, $false, $false, $true, 1, $false, "35×86=", 2) | Out-Null
$d.Content.Find.Execute("59×69=", $true, $false, $false, $false, $false, $true, 1, $false, "68×18=", 2) | Out-Null
$d.Content.Find.Execute("41×38=", $true, $false, $false, $false, $false, $true, 1, $false, "49×65=", 2) | Out-Null
$d.Content.Find.Execute("74×53=", $true, $false, $false, $false, $false, $true, 1, $false, "88×63=", 2) | Out-Null
$d.Content.Find.Execute("33×59=", $true, $false, $false, $false, $false, $true, 1, $false, "97×65=", 2) | Out-Null
$d.Content.Find.Execute("28×75=", $true, $false, $false, $false, $false, $true, 1, $false, "33×30=", 2) | Out-Null
$d.Content.Find.Execute("43×71=", $true, $false, $false, $false, $false, $true, 1, $false, "98×81=", 2) | Out-Null
$d.Content.Find.Execute("90×30=", $true, $false, $false, $false, $false, $true, 1, $false, "11×97=", 2) | Out-Null
$d.Content.Find.Execute("53×98=", $true, $false, $false, $false, $false, $true, 1, $false, "98×90=", 2) | Out-Null
$d.Content.Find.Execute("11×35=", $true, $false, $false, $false, $false, $true, 1, $false, "92×94=", 2) | Out-Null
$d.Content.Find.Execute("23×22=", $true, $false, $false, $false, $false, $true, 1, $false, "28×49=", 2) | Out-Null
$d.Content.Find.Execute("20×37=", $true, $false, $false, $false, $false, $true, 1, $false, "54×61=", 2) | Out-Null
$d.Content.Find.Execute("82×40=", $true, $false, $false, $false, $false, $true, 1, $false, "17×38=", 2) | Out-Null
$d.Content.Find.Execute("85×84=", $true, $false, $false, $false, $false, $true, 1, $false, "75×58=", 2) | Out-Null
$d.Content.Find.Execute("25×74=", $true, $false, $false, $false, $false, $true, 1, $false, "24×53=", 2) | Out-Null
$d.Content.Find.Execute("63×27=", $true, $false, $false, $false, $false, $true, 1, $false, "34×77=", 2) | Out-Null
$d.Content.Find.Execute("57×87=", $true, $false, $false, $false, $false, $true, 1, $false, "24×42=", 2) | Out-Null
$d.Content.Find.Execute("87×20=", $true, $false, $false, $false, $false, $true, 1, $false, "37×33=", 2) | Out-Null
$d.Content.Find.Execute("72×27=", $true, $false, $false, $false, $false, $true, 1, $false, "14×75=", 2) | Out-Null
$d.Content.Find.Execute("18×56=", $true, $false, $false, $false, $false, $true, 1, $false, "50×20=", 2) | Out-Null
